# Applies the 13-Abril-2020 13:22 COVID data refresh to the "Pais" sheet.
# Several countries (Bielorrusia, Afganistan, Nepal, ...) overtook their
# neighbours in the total-cases ranking, so the rows below them shift down
# by one position while keeping their own totals unless explicitly listed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 13 de Abril de 2020 a las 13:22"

# Row 11: Iran
$ws.Cells.Item(11, 1).Value = "Iran"
$ws.Cells.Item(11, 2).Value = 73303
$ws.Cells.Item(11, 3).Value = 1617
$ws.Cells.Item(11, 4).Value = 45983
$ws.Cells.Item(11, 5).Value = 22735
$ws.Cells.Item(11, 6).Value = 3877
$ws.Cells.Item(11, 7).Value = 111
$ws.Cells.Item(11, 8).Value = 4585

# Row 20: Austria
$ws.Cells.Item(20, 1).Value = "Austria"
$ws.Cells.Item(20, 2).Value = 13984
$ws.Cells.Item(20, 3).Value = 39
$ws.Cells.Item(20, 4).Value = 7343
$ws.Cells.Item(20, 5).Value = 6273
$ws.Cells.Item(20, 6).Value = 239
$ws.Cells.Item(20, 7).Value = 18
$ws.Cells.Item(20, 8).Value = 368

# Row 48: Finlandia
$ws.Cells.Item(48, 1).Value = "Finlandia"
$ws.Cells.Item(48, 2).Value = 3064
$ws.Cells.Item(48, 3).Value = 90
$ws.Cells.Item(48, 4).Value = 300
$ws.Cells.Item(48, 5).Value = 2705
$ws.Cells.Item(48, 6).Value = 74
$ws.Cells.Item(48, 7).Value = 3
$ws.Cells.Item(48, 8).Value = 59

# Row 50: Bielorrusia
$ws.Cells.Item(50, 1).Value = "Bielorrusia"
$ws.Cells.Item(50, 2).Value = 2919
$ws.Cells.Item(50, 3).Value = 341
$ws.Cells.Item(50, 4).Value = 203
$ws.Cells.Item(50, 5).Value = 2687
$ws.Cells.Item(50, 6).Value = 55
$ws.Cells.Item(50, 7).Value = 3
$ws.Cells.Item(50, 8).Value = 29

# Row 51: Colombia
$ws.Cells.Item(51, 1).Value = "Colombia"
$ws.Cells.Item(51, 2).Value = 2776
$ws.Cells.Item(51, 3).Value = 0
$ws.Cells.Item(51, 4).Value = 270
$ws.Cells.Item(51, 5).Value = 2397
$ws.Cells.Item(51, 6).Value = 92
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 109

# Row 52: Tailandia
$ws.Cells.Item(52, 1).Value = "Tailandia"
$ws.Cells.Item(52, 2).Value = 2579
$ws.Cells.Item(52, 3).Value = 28
$ws.Cells.Item(52, 4).Value = 1288
$ws.Cells.Item(52, 5).Value = 1251
$ws.Cells.Item(52, 6).Value = 61
$ws.Cells.Item(52, 7).Value = 2
$ws.Cells.Item(52, 8).Value = 40

# Row 73: Bosnia y Herzegovina
$ws.Cells.Item(73, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(73, 2).Value = 1030
$ws.Cells.Item(73, 3).Value = 21
$ws.Cells.Item(73, 4).Value = 206
$ws.Cells.Item(73, 5).Value = 785
$ws.Cells.Item(73, 6).Value = 4
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 39

# Row 76: Uzbekistan
$ws.Cells.Item(76, 1).Value = "Uzbekistan"
$ws.Cells.Item(76, 2).Value = 896
$ws.Cells.Item(76, 3).Value = 31
$ws.Cells.Item(76, 4).Value = 68
$ws.Cells.Item(76, 5).Value = 824
$ws.Cells.Item(76, 6).Value = 8
$ws.Cells.Item(76, 7).Value = 0
$ws.Cells.Item(76, 8).Value = 4

# Row 86: Afganistan
$ws.Cells.Item(86, 1).Value = "Afganistan"
$ws.Cells.Item(86, 2).Value = 665
$ws.Cells.Item(86, 3).Value = 58
$ws.Cells.Item(86, 4).Value = 32
$ws.Cells.Item(86, 5).Value = 612
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 3
$ws.Cells.Item(86, 8).Value = 21

# Row 87: Letonia
$ws.Cells.Item(87, 1).Value = "Letonia"
$ws.Cells.Item(87, 2).Value = 653
$ws.Cells.Item(87, 3).Value = 2
$ws.Cells.Item(87, 4).Value = 16
$ws.Cells.Item(87, 5).Value = 632
$ws.Cells.Item(87, 6).Value = 2
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 5

# Row 88: Principado de Andorra
$ws.Cells.Item(88, 1).Value = "Principado de Andorra"
$ws.Cells.Item(88, 2).Value = 638
$ws.Cells.Item(88, 3).Value = 0
$ws.Cells.Item(88, 4).Value = 128
$ws.Cells.Item(88, 5).Value = 481
$ws.Cells.Item(88, 6).Value = 17
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 29

# Row 89: Republica de Chipre
$ws.Cells.Item(89, 1).Value = "Republica de Chipre"
$ws.Cells.Item(89, 2).Value = 633
$ws.Cells.Item(89, 3).Value = 0
$ws.Cells.Item(89, 4).Value = 65
$ws.Cells.Item(89, 5).Value = 557
$ws.Cells.Item(89, 6).Value = 8
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 8).Value = 11

# Row 90: Libano
$ws.Cells.Item(90, 1).Value = "Libano"
$ws.Cells.Item(90, 2).Value = 632
$ws.Cells.Item(90, 3).Value = 2
$ws.Cells.Item(90, 4).Value = 80
$ws.Cells.Item(90, 5).Value = 532
$ws.Cells.Item(90, 6).Value = 34
$ws.Cells.Item(90, 7).Value = 0
$ws.Cells.Item(90, 8).Value = 20

# Row 112: Vietnam
$ws.Cells.Item(112, 1).Value = "Vietnam"
$ws.Cells.Item(112, 2).Value = 265
$ws.Cells.Item(112, 3).Value = 3
$ws.Cells.Item(112, 4).Value = 145
$ws.Cells.Item(112, 5).Value = 120
$ws.Cells.Item(112, 6).Value = 8
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 0

# Row 182: Nepal
$ws.Cells.Item(182, 1).Value = "Nepal"
$ws.Cells.Item(182, 2).Value = 14
$ws.Cells.Item(182, 3).Value = 2
$ws.Cells.Item(182, 4).Value = 1
$ws.Cells.Item(182, 5).Value = 13
$ws.Cells.Item(182, 6).Value = 0
$ws.Cells.Item(182, 7).Value = 0
$ws.Cells.Item(182, 8).Value = 0

# Row 183: Belice
$ws.Cells.Item(183, 1).Value = "Belice"
$ws.Cells.Item(183, 2).Value = 14
$ws.Cells.Item(183, 3).Value = 0
$ws.Cells.Item(183, 4).Value = 0
$ws.Cells.Item(183, 5).Value = 12
$ws.Cells.Item(183, 6).Value = 1
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 2

# Row 184: Zimbabue
$ws.Cells.Item(184, 1).Value = "Zimbabue"
$ws.Cells.Item(184, 2).Value = 14
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(184, 4).Value = 0
$ws.Cells.Item(184, 5).Value = 11
$ws.Cells.Item(184, 6).Value = 0
$ws.Cells.Item(184, 7).Value = 0
$ws.Cells.Item(184, 8).Value = 3

# Row 185: Suazilandia
$ws.Cells.Item(185, 1).Value = "Suazilandia"
$ws.Cells.Item(185, 2).Value = 14
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(185, 4).Value = 7
$ws.Cells.Item(185, 5).Value = 7
$ws.Cells.Item(185, 6).Value = 0
$ws.Cells.Item(185, 7).Value = 0
$ws.Cells.Item(185, 8).Value = 0

# Row 186: Curazao
$ws.Cells.Item(186, 1).Value = "Curazao"
$ws.Cells.Item(186, 2).Value = 14
$ws.Cells.Item(186, 3).Value = 0
$ws.Cells.Item(186, 4).Value = 8
$ws.Cells.Item(186, 5).Value = 5
$ws.Cells.Item(186, 6).Value = 0
$ws.Cells.Item(186, 7).Value = 0
$ws.Cells.Item(186, 8).Value = 1

